$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "AEDB.CEA"
$ws.Cells.Item(2, 2).Value = "epmajor.30days"
$ws.Cells.Item(2, 3).Value = "MCP1_pg_ug_2015_rank"
$ws.Cells.Item(2, 4).Value = 1.11962496853713
$ws.Cells.Item(2, 5).Value = 0.359458849191985
$ws.Cells.Item(2, 6).Value = 3.06370500204444
$ws.Cells.Item(2, 7).Value = 1.5145004129068
$ws.Cells.Item(2, 8).Value = 6.19761358898339
$ws.Cells.Item(2, 9).Value = 3.11475144110069
$ws.Cells.Item(2, 10).Value = 0.00184099973150297
$ws.Cells.Item(2, 11).Value = 1186
$ws.Cells.Item(2, 12).Value = 40

# Row 3
$ws.Cells.Item(3, 1).Value = "AEDB.CEA"
$ws.Cells.Item(3, 2).Value = "epmajor.30days"
$ws.Cells.Item(3, 3).Value = "MCP1_pg_ml_2015_rank"
$ws.Cells.Item(3, 4).Value = 0.820847048319384
$ws.Cells.Item(3, 5).Value = 0.358633874094218
$ws.Cells.Item(3, 6).Value = 2.27242387536542
$ws.Cells.Item(3, 7).Value = 1.12515933315019
$ws.Cells.Item(3, 8).Value = 4.589492454258
$ws.Cells.Item(3, 9).Value = 2.28881627646734
$ws.Cells.Item(3, 10).Value = 0.022090029371139
$ws.Cells.Item(3, 11).Value = 1187
$ws.Cells.Item(3, 12).Value = 40

# Row 4
$ws.Cells.Item(4, 1).Value = "AEDB.CEA"
$ws.Cells.Item(4, 2).Value = "epmajor.30days"
$ws.Cells.Item(4, 3).Value = "MCP1_rank"
$ws.Cells.Item(4, 4).Value = -0.0533758162933218
$ws.Cells.Item(4, 5).Value = 0.504576631738758
$ws.Cells.Item(4, 6).Value = 0.948023662786552
$ws.Cells.Item(4, 7).Value = 0.35262644405829
$ws.Cells.Item(4, 8).Value = 2.54872792539253
$ws.Cells.Item(4, 9).Value = -0.105783369533762
$ws.Cells.Item(4, 10).Value = 0.915754231936841
$ws.Cells.Item(4, 11).Value = 549
$ws.Cells.Item(4, 12).Value = 17

# Row 5
$ws.Cells.Item(5, 1).Value = "AEDB.CEA"
$ws.Cells.Item(5, 2).Value = "epstroke.30days"
$ws.Cells.Item(5, 3).Value = "MCP1_pg_ug_2015_rank"
$ws.Cells.Item(5, 4).Value = 0.907885502034972
$ws.Cells.Item(5, 5).Value = 0.377280671418516
$ws.Cells.Item(5, 6).Value = 2.47907498836033
$ws.Cells.Item(5, 7).Value = 1.18342808808882
$ws.Cells.Item(5, 8).Value = 5.193228773063
$ws.Cells.Item(5, 9).Value = 2.40639282850475
$ws.Cells.Item(5, 10).Value = 0.0161109303441679
$ws.Cells.Item(5, 11).Value = 1186
$ws.Cells.Item(5, 12).Value = 34

# Row 6
$ws.Cells.Item(6, 1).Value = "AEDB.CEA"
$ws.Cells.Item(6, 2).Value = "epstroke.30days"
$ws.Cells.Item(6, 3).Value = "MCP1_pg_ml_2015_rank"
$ws.Cells.Item(6, 4).Value = 0.733671201916969
$ws.Cells.Item(6, 5).Value = 0.383138843650386
$ws.Cells.Item(6, 6).Value = 2.08271264838825
$ws.Cells.Item(6, 7).Value = 0.982867527261892
$ws.Cells.Item(6, 8).Value = 4.41330276506389
$ws.Cells.Item(6, 9).Value = 1.91489642482307
$ws.Cells.Item(6, 10).Value = 0.0555057138602618
$ws.Cells.Item(6, 11).Value = 1187
$ws.Cells.Item(6, 12).Value = 34

# Row 7
$ws.Cells.Item(7, 1).Value = "AEDB.CEA"
$ws.Cells.Item(7, 2).Value = "epstroke.30days"
$ws.Cells.Item(7, 3).Value = "MCP1_rank"
$ws.Cells.Item(7, 4).Value = 0.208084795502029
$ws.Cells.Item(7, 5).Value = 0.550356875897256
$ws.Cells.Item(7, 6).Value = 1.23131757531419
$ws.Cells.Item(7, 7).Value = 0.418694106713552
$ws.Cells.Item(7, 8).Value = 3.62112326628682
$ws.Cells.Item(7, 9).Value = 0.378090661923294
$ws.Cells.Item(7, 10).Value = 0.705363245048678
$ws.Cells.Item(7, 11).Value = 549
$ws.Cells.Item(7, 12).Value = 14

# Row 8
$ws.Cells.Item(8, 1).Value = "AEDB.CEA"
$ws.Cells.Item(8, 2).Value = "epcoronary.30days"
$ws.Cells.Item(8, 3).Value = "MCP1_pg_ug_2015_rank"
$ws.Cells.Item(8, 4).Value = 0.987619627716873
$ws.Cells.Item(8, 5).Value = 0.704674023443691
$ws.Cells.Item(8, 6).Value = 2.68483595085675
$ws.Cells.Item(8, 7).Value = 0.674663346647003
$ws.Cells.Item(8, 8).Value = 10.6843570483523
$ws.Cells.Item(8, 9).Value = 1.40152693991819
$ws.Cells.Item(8, 10).Value = 0.161056557343573
$ws.Cells.Item(8, 11).Value = 1186
$ws.Cells.Item(8, 12).Value = 11

# Row 9
$ws.Cells.Item(9, 1).Value = "AEDB.CEA"
$ws.Cells.Item(9, 2).Value = "epcoronary.30days"
$ws.Cells.Item(9, 3).Value = "MCP1_pg_ml_2015_rank"
$ws.Cells.Item(9, 4).Value = 0.931359355801676
$ws.Cells.Item(9, 5).Value = 0.725473263558908
$ws.Cells.Item(9, 6).Value = 2.53795682012838
$ws.Cells.Item(9, 7).Value = 0.612278373777202
$ws.Cells.Item(9, 8).Value = 10.5200919985131
$ws.Cells.Item(9, 9).Value = 1.28379556157971
$ws.Cells.Item(9, 10).Value = 0.199213496991895
$ws.Cells.Item(9, 11).Value = 1187
$ws.Cells.Item(9, 12).Value = 11

# Row 10
$ws.Cells.Item(10, 1).Value = "AEDB.CEA"
$ws.Cells.Item(10, 2).Value = "epcoronary.30days"
$ws.Cells.Item(10, 3).Value = "MCP1_rank"
$ws.Cells.Item(10, 4).Value = -0.0472681771167935
$ws.Cells.Item(10, 5).Value = 0.933970603742045
$ws.Cells.Item(10, 6).Value = 0.953831567488339
$ws.Cells.Item(10, 7).Value = 0.15291844170768
$ws.Cells.Item(10, 8).Value = 5.94954178827189
$ws.Cells.Item(10, 9).Value = -0.0506099195492973
$ws.Cells.Item(10, 10).Value = 0.959636358299397
$ws.Cells.Item(10, 11).Value = 549
$ws.Cells.Item(10, 12).Value = 5

# Row 11
$ws.Cells.Item(11, 1).Value = "AEDB.CEA"
$ws.Cells.Item(11, 2).Value = "epcvdeath.30days"
$ws.Cells.Item(11, 3).Value = "MCP1_pg_ug_2015_rank"
$ws.Cells.Item(11, 4).Value = 1.62403826178446
$ws.Cells.Item(11, 5).Value = 1.26373486532014
$ws.Cells.Item(11, 6).Value = 5.07353727538219
$ws.Cells.Item(11, 7).Value = 0.426184863807573
$ws.Cells.Item(11, 8).Value = 60.3981573975249
$ws.Cells.Item(11, 9).Value = 1.28510995965363
$ws.Cells.Item(11, 10).Value = 0.198753861692165
$ws.Cells.Item(11, 11).Value = 1186
$ws.Cells.Item(11, 12).Value = 4

# Row 12
$ws.Cells.Item(12, 1).Value = "AEDB.CEA"
$ws.Cells.Item(12, 2).Value = "epcvdeath.30days"
$ws.Cells.Item(12, 3).Value = "MCP1_pg_ml_2015_rank"
$ws.Cells.Item(12, 4).Value = -0.255589970332927
$ws.Cells.Item(12, 5).Value = 1.15243760621838
$ws.Cells.Item(12, 6).Value = 0.774459455050924
$ws.Cells.Item(12, 7).Value = 0.0809140597320861
$ws.Cells.Item(12, 8).Value = 7.41264805527896
$ws.Cells.Item(12, 9).Value = -0.221782046120157
$ws.Cells.Item(12, 10).Value = 0.824483556544668
$ws.Cells.Item(12, 11).Value = 1187
$ws.Cells.Item(12, 12).Value = 4

# Row 13
$ws.Cells.Item(13, 1).Value = "AEDB.CEA"
$ws.Cells.Item(13, 2).Value = "epcvdeath.30days"
$ws.Cells.Item(13, 3).Value = "MCP1_rank"
$ws.Cells.Item(13, 4).Value = 140.863868031622
$ws.Cells.Item(13, 5).Value = 34819.2952907117
$ws.Cells.Item(13, 6).Value = "15010687568824498783824533463539265762500444067615544447074304"
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = "#NUM!"
$ws.Cells.Item(13, 9).Value = 0.00404556918385419
$ws.Cells.Item(13, 10).Value = 0.996772111613488
$ws.Cells.Item(13, 11).Value = 549
$ws.Cells.Item(13, 12).Value = 1
